$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the existing "insertQuery()" sheet (placed right after it) so the
# new sheet inherits the same layout, column widths, merged cells and styles.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "insertFile()"

# --- Update the headings/labels that are specific to insertFile() ---
$ws2.Range("A2").Value = "EQUIVALENCE CLASS PARTITIONING insertFile()"
$ws2.Range("A8").Value = "BOUNDARY VALUE ANALYSIS insertFile()"
$ws2.Range("A12").Value = "filePath"
$ws2.Range("A15").Value = "TEST CASES insertFile()"
$ws2.Range("D17").Value = "filePath"

# --- Row 18: "Invalid Path" / fm (FileManager) being declared test case ---
$ws2.Range("B18").Value = "Invalid Path"
$ws2.Range("C18").Value = "fm (FileManager) being declared"
$ws2.Range("D18").Value = "Invalid path"
$ws2.Range("E18").Value = "False "
$ws2.Range("F18").Value = "False "
$ws2.Range("G18").Value = "Passed"

# --- Row 19: "Invalid Path" / Null path test case ---
$ws2.Range("B19").Value = "Invalid Path"
$ws2.Range("D19").Value = "Null path"
$ws2.Range("E19").Value = "False "
$ws2.Range("F19").Value = "False "
$ws2.Range("G19").Value = "Passed"

# --- Row 20 (new): "Valid Path" / Valid path test case ---
$ws2.Range("A20").Value = 2
$ws2.Range("B20").Value = "Valid Path"
$ws2.Range("D20").Value = "Valid path"
$ws2.Range("E20").Value = "True "
$ws2.Range("F20").Value = "True "
$ws2.Range("G20").Value = "Passed"

# Style A20/B20/D20/E20/F20/G20 like the row above (row 19) and C20 like the
# plain bordered cell used elsewhere in that merged column, so the new row
# matches the look of the rest of the test-case table.
$ws2.Range("A20").Style = $ws2.Range("A19").Style
$ws2.Range("B20").Style = $ws2.Range("B19").Style
$ws2.Range("C20").Style = $ws2.Range("C17").Style
$ws2.Range("D20").Style = $ws2.Range("D19").Style
$ws2.Range("E20").Style = $ws2.Range("E19").Style
$ws2.Range("F20").Style = $ws2.Range("F19").Style
$ws2.Range("G20").Style = $ws2.Range("G19").Style

# Extend the merged "fm (FileManager) being declared" cell down to cover the
# new row (C18:C19 -> C18:C20), matching the target layout.
$ws2.Range("C18:C19").UnMerge()
$ws2.Range("C18:C20").Merge()

$ws1.Select()
